# pure tension calculation complete
# - membrane_density switched from g/cm^3 (1.2) to kg/m^3 (1400)
# - new "Results" sheet added after "input" with Tension_T / RimForce outputs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update membrane_density row (row 4) on the "input" sheet ---
$ws.Range("B4").Value = 1400
$ws.Range("C4").Value = "kg/m^3"

# --- add the "Results" sheet right after "input" ---
$results = $wb.Worksheets.Add($null, $ws)
$results.Name = "Results"

$results.Range("A1").Value = "Tension_T"
$results.Range("B1").Value = 887.891240031866
$results.Range("C1").Value = "kN/m"

$results.Range("A2").Value = "RimForce"
$results.Range("B2").Value = 17004.13727052465
$results.Range("C2").Value = "kN"

# approximate the autofit column widths captured in the saved workbook
$results.Columns.Item(1).ColumnWidth = 10.0703125
$results.Columns.Item(2).ColumnWidth = 9.92578125
$results.Columns.Item(3).ColumnWidth = 5.78515625

# go back to the "input" sheet and leave the selection on A5, matching the saved state
$ws.Activate() | Out-Null
$ws.Range("A5").Select() | Out-Null

# force full recalculation next time the workbook is opened
$wb.Application.CalculateFullRebuild() | Out-Null
